$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.554.04'
$ws.Range("E2").Value = '  +7.92%  '
$ws.Range("D3").Value = '3.644.43'
$ws.Range("E3").Value = '  +4.76%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '419.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.664'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.50%  '
$ws.Range("D8").Value = '3.626.44'
$ws.Range("E8").Value = '  +4.49%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.202'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +33.43%  '
$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.776'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000443'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +101.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '4.222.75'
$ws.Range("E15").Value = '  +4.74%  '
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.652.75'
$ws.Range("E17").Value = '  +4.76%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("E19").Value = '  +3.09%  '
$ws.Range("D20").Value = '68.547.26'
$ws.Range("E20").Value = '  +8.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '462.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '90.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '35.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.06%  '
$ws.Range("E30").Value = '  +5.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.119'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.69%  '
$ws.Range("E34").Value = '  -4.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '40.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").Value = '0.0₃0838'
$ws.Range("E36").Value = '  +27.38%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0487'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.149'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '148.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.43%  '
$ws.Range("E46").Value = '  -7.03%  '
$ws.Range("E47").Value = '  +21.43%  '
$ws.Range("B48").Value = 'TheGraph'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.305'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.24%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.83%  '
$ws.Range("E50").Value = '  -1.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.67'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.67%  '
